$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "dMubi476"
$ws.Range("B2").Value = 231106318
$ws.Range("C2").Value = "xmnhxnm60"
$ws.Range("D2").Value = "HX&!43wv"
$ws.Range("F2").Value = "vItZBynr"
$ws.Range("G2").Value = "SWms"

# Row 3
$ws.Range("A3").Value = "OWpYN334"
$ws.Range("B3").Value = 231106317
$ws.Range("C3").Value = "bykekoz48"
$ws.Range("D3").Value = "FGd5%k6#"
$ws.Range("F3").Value = "KyFKAOPa"
$ws.Range("G3").Value = "shAo"
